# Grado06 "Seguimiento" sheet update:
#  - CS_04_09_CO: recepción (B6) date corrected
#  - CS_06_04_CO: recepción dates entered for items 2 and 4 (B7, B9)
#  - leave the cursor/selection where the user was last working (E9)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Seguimiento")

# Row 6 (item 1): recepción date corrected 3/23/2015 -> 2/23/2015
$ws.Range("B6").Value = 42058

# Row 7 (item 2): recepción date entered 3/9/2015
$ws.Range("B7").Value = 42072

# Row 9 (item 4): recepción date entered 3/15/2015
$ws.Range("B9").Value = 42078

# Update the active selection to reflect where work left off
$null = $ws.Range("E9").Select()
